$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H52").Value = 8335
$ws.Range("J52").Value = 8335
$ws.Range("L52").Value = 25005
$ws.Range("N52").Value = -25325
$ws.Range("H62").Value = 3252.1082
$ws.Range("I62").Value = 2766.0476
$ws.Range("J62").Value = 3890.0625
$ws.Range("K62").Value = 2766.0476
$ws.Range("L62").Value = 3890.0625
$ws.Range("M62").Value = -2142.0476
$ws.Range("N62").Value = -5138.0625
$ws.Range("H65").Value = 3252.1082
$ws.Range("I65").Value = 2766.0476
$ws.Range("J65").Value = 3890.0625
$ws.Range("K65").Value = 13830.238
$ws.Range("L65").Value = 19450.3125
$ws.Range("M65").Value = -10710.238
$ws.Range("N65").Value = -25690.3125
$ws.Range("H86").Value = 7224
$ws.Range("I86").Value = 677.3
$ws.Range("J86").Value = 20317.4
$ws.Range("K86").Value = 677.3
$ws.Range("L86").Value = 20317.4
$ws.Range("M86").Value = 445.7
$ws.Range("N86").Value = -22563.4
$ws.Range("H89").Value = 7224
$ws.Range("I89").Value = 677.3
$ws.Range("J89").Value = 20317.4
$ws.Range("K89").Value = 3386.5
$ws.Range("L89").Value = 101587
$ws.Range("M89").Value = 2229.5
$ws.Range("N89").Value = -112819
$ws.Range("H92").Value = 58823970
$ws.Range("I92").Value = 100000424
$ws.Range("J92").Value = 456.85715
$ws.Range("K92").Value = 100000424
$ws.Range("L92").Value = 456.85715
$ws.Range("M92").Value = -99999176
$ws.Range("N92").Value = -2952.85715
$ws.Range("H96").Value = 20833688
$ws.Range("I96").Value = 22727590
$ws.Range("J96").Value = 776
$ws.Range("K96").Value = 68182770
$ws.Range("L96").Value = 2328
$ws.Range("M96").Value = -68181397
$ws.Range("N96").Value = -5074
$ws.Range("H112").Value = 4116313.8
$ws.Range("J112").Value = 4116313.8
$ws.Range("L112").Value = 12348941.4
$ws.Range("N112").Value = -12351157.4
$ws.Range("H129").Value = 245055.19
$ws.Range("J129").Value = 264379.28
$ws.Range("L129").Value = 793137.8400000001
$ws.Range("N129").Value = -803137.8400000001
$ws.Range("H132").Value = 2116.5833
$ws.Range("I132").Value = 2207.3
$ws.Range("J132").Value = 1663
$ws.Range("K132").Value = 6621.900000000001
$ws.Range("L132").Value = 4989
$ws.Range("M132").Value = -4091.900000000001
$ws.Range("N132").Value = -10049
$ws.Range("H137").Value = 1977.2759
$ws.Range("I137").Value = 1981.3684
$ws.Range("J137").Value = 1969.5
$ws.Range("K137").Value = 5944.1052
$ws.Range("L137").Value = 5908.5
$ws.Range("M137").Value = -3394.1052
$ws.Range("N137").Value = -11008.5
$ws.Range("H138").Value = 10528455
$ws.Range("I138").Value = 20408928
$ws.Range("J138").Value = 3601.9348
$ws.Range("K138").Value = 61226784
$ws.Range("L138").Value = 10805.8044
$ws.Range("M138").Value = -61221644
$ws.Range("N138").Value = -21085.8044
$ws.Range("H141").Value = 1036.4231
$ws.Range("I141").Value = 783.125
$ws.Range("K141").Value = 2349.375
$ws.Range("M141").Value = 2830.625

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2656.02
$ws.Range("I32").Value = 2397.6667
$ws.Range("J32").Value = 4981.2
$ws.Range("K32").Value = 2397.6667
$ws.Range("L32").Value = 4981.2
$ws.Range("M32").Value = -2110.6667
$ws.Range("N32").Value = -5555.2
$ws.Range("H53").Value = 8349.333000000001
$ws.Range("I53").Value = 5024
$ws.Range("J53").Value = 15000
$ws.Range("K53").Value = 5024
$ws.Range("L53").Value = 15000
$ws.Range("M53").Value = -4342
$ws.Range("N53").Value = -16364
$ws.Range("H63").Value = 1801.3334
$ws.Range("J63").Value = 1199
$ws.Range("L63").Value = 1199
$ws.Range("N63").Value = -2571
$ws.Range("H66").Value = 1801.3334
$ws.Range("J66").Value = 1199
$ws.Range("L66").Value = 5995
$ws.Range("N66").Value = -12859
$ws.Range("H74").Value = 2480.0688
$ws.Range("I74").Value = 2404.96
$ws.Range("K74").Value = 2404.96
$ws.Range("M74").Value = -1530.96
$ws.Range("H77").Value = 2480.0688
$ws.Range("I77").Value = 2404.96
$ws.Range("K77").Value = 12024.8
$ws.Range("M77").Value = -7656.799999999999
$ws.Range("H97").Value = 1279.9714
$ws.Range("I97").Value = 1128.1111
$ws.Range("K97").Value = 1128.1111
$ws.Range("M97").Value = -632.1111000000001
$ws.Range("H132").Value = 13452.262
$ws.Range("I132").Value = 1502.25
$ws.Range("J132").Value = 85152.336
$ws.Range("K132").Value = 4506.75
$ws.Range("L132").Value = 255457.008
$ws.Range("M132").Value = -1976.75
$ws.Range("N132").Value = -260517.008

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1229.4445
$ws.Range("I94").Value = 772.1429000000001
$ws.Range("J94").Value = 2830
$ws.Range("K94").Value = 772.1429000000001
$ws.Range("L94").Value = 2830
$ws.Range("M94").Value = -321.1429000000001
$ws.Range("N94").Value = -3732
$ws.Range("H99").Value = 1755.4546
$ws.Range("I99").Value = 1902.5
$ws.Range("K99").Value = 1902.5
$ws.Range("M99").Value = -404.5
$ws.Range("H134").Value = 2707.2075
$ws.Range("I134").Value = 2928.9524
$ws.Range("K134").Value = 8786.8572
$ws.Range("M134").Value = -6251.8572

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2602.717
$ws.Range("I31").Value = 1496.2821
$ws.Range("K31").Value = 1496.2821
$ws.Range("M31").Value = -1201.2821
$ws.Range("H34").Value = 2602.717
$ws.Range("I34").Value = 1496.2821
$ws.Range("K34").Value = 1496.2821
$ws.Range("M34").Value = -1294.2821
$ws.Range("H58").Value = 15815.941
$ws.Range("I58").Value = 1043.2142
$ws.Range("J58").Value = 84755.336
$ws.Range("K58").Value = 1043.2142
$ws.Range("L58").Value = 84755.336
$ws.Range("M58").Value = -840.2141999999999
$ws.Range("N58").Value = -85161.336
$ws.Range("H99").Value = 23813410
$ws.Range("I99").Value = 3722.2222
$ws.Range("J99").Value = 41670676
$ws.Range("K99").Value = 3722.2222
$ws.Range("L99").Value = 41670676
$ws.Range("M99").Value = -2224.2222
$ws.Range("N99").Value = -41673672
$ws.Range("H126").Value = 23813410
$ws.Range("I126").Value = 3722.2222
$ws.Range("J126").Value = 41670676
$ws.Range("K126").Value = 11166.6666
$ws.Range("L126").Value = 125012028
$ws.Range("M126").Value = -8696.6666
$ws.Range("N126").Value = -125016968
$ws.Range("H132").Value = 1704.3922
$ws.Range("I132").Value = 1416.0454
$ws.Range("J132").Value = 3516.8572
$ws.Range("K132").Value = 4248.1362
$ws.Range("L132").Value = 10550.5716
$ws.Range("M132").Value = -1718.1362
$ws.Range("N132").Value = -15610.5716
$ws.Range("H134").Value = 884.7174
$ws.Range("I134").Value = 787.7560999999999
$ws.Range("J134").Value = 1679.8
$ws.Range("K134").Value = 2363.2683
$ws.Range("L134").Value = 5039.4
$ws.Range("M134").Value = 171.7317000000003
$ws.Range("N134").Value = -10109.4
$ws.Range("H136").Value = 15815.941
$ws.Range("I136").Value = 1043.2142
$ws.Range("J136").Value = 84755.336
$ws.Range("K136").Value = 3129.6426
$ws.Range("L136").Value = 254266.008
$ws.Range("M136").Value = -579.6425999999997
$ws.Range("N136").Value = -259366.008

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1420.0454
$ws.Range("I5").Value = 1032.9
$ws.Range("J5").Value = 1742.6666
$ws.Range("K5").Value = 3098.7
$ws.Range("L5").Value = 5227.9998
$ws.Range("M5").Value = -2986.7
$ws.Range("N5").Value = -5451.9998
$ws.Range("H20").Value = 1766.6666
$ws.Range("I20").Value = 1766.6666
$ws.Range("K20").Value = 5299.9998
$ws.Range("M20").Value = -5072.9998
$ws.Range("H131").Value = 730.91
$ws.Range("J131").Value = 743.84045
$ws.Range("L131").Value = 2231.52135
$ws.Range("N131").Value = -12311.52135
$ws.Range("H135").Value = 1420.0454
$ws.Range("I135").Value = 1032.9
$ws.Range("J135").Value = 1742.6666
$ws.Range("K135").Value = 9296.1
$ws.Range("L135").Value = 15683.9994
$ws.Range("M135").Value = -6761.1
$ws.Range("N135").Value = -20753.9994

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1228.8948
$ws.Range("I97").Value = 1241.6111
$ws.Range("K97").Value = 1241.6111
$ws.Range("M97").Value = -745.6111000000001
$ws.Range("H132").Value = 21244.963
$ws.Range("I132").Value = 2900.5833
$ws.Range("J132").Value = 168000
$ws.Range("K132").Value = 8701.749899999999
$ws.Range("L132").Value = 504000
$ws.Range("M132").Value = -6171.749899999999
$ws.Range("N132").Value = -509060

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2280.1333
$ws.Range("I82").Value = 2045.6364
$ws.Range("K82").Value = 2045.6364
$ws.Range("M82").Value = -1684.6364
$ws.Range("H85").Value = 2280.1333
$ws.Range("I85").Value = 2045.6364
$ws.Range("K85").Value = 2045.6364
$ws.Range("M85").Value = -797.6364000000001
$ws.Range("H93").Value = 961.64703
$ws.Range("J93").Value = 963.4
$ws.Range("L93").Value = 963.4
$ws.Range("N93").Value = -3459.4
$ws.Range("H100").Value = 2213.4285
$ws.Range("I100").Value = 1483.6666
$ws.Range("J100").Value = 2760.75
$ws.Range("K100").Value = 1483.6666
$ws.Range("L100").Value = 2760.75
$ws.Range("M100").Value = -942.6666
$ws.Range("N100").Value = -3842.75
$ws.Range("H132").Value = 1190.8
$ws.Range("I132").Value = 1267.5927
$ws.Range("K132").Value = 3802.7781
$ws.Range("M132").Value = -1272.7781

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 538
$ws.Range("I100").Value = 597.5
$ws.Range("K100").Value = 1195
$ws.Range("M100").Value = -654
$ws.Range("H136").Value = 15874949
$ws.Range("I136").Value = 22728186
$ws.Range("K136").Value = 68184558
$ws.Range("M136").Value = -68182008
